{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === \"Docente(s) Respons\u00e1vel(eis)\") {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Target paragraph not found\");\n}\n\nconst newPara = target.insertParagraph(\"6712818 - Mauricio Lamano Ferreira\", Word.InsertLocation.after);\nnewPara.style = \"List Bullet\";\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the \"Docente(s) Respons\u00e1vel(eis)\" heading paragraph.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Trim() -eq \"Docente(s) Respons\u00e1vel(eis)\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find the 'Docente(s) Respons\u00e1vel(eis)' paragraph\"\n}\n\n# Insert a new paragraph right after it and fill it in as a bulleted\n# list item with the instructor's name.\n$target.Range.InsertParagraphAfter()\n$newPara = $target.Next()\n$newPara.Range.Text = \"6712818 - Mauricio Lamano Ferreira\"\n$newPara.Style = $d.Styles.Item(\"List Bullet\")\n\n$d.Save()\n"}
